{"js": "// --- Change 1: title \"JUSTIFICACION\" -> \"JUSTIFICACI\u00d3N\" -------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst titleRange = titlePara.getRange();\n// Retype the title with the missing accent fixed, keeping the run's\n// formatting (Arial, bold, size 28) which insertText(Replace) preserves.\ntitleRange.insertText(\"JUSTIFICACI\u00d3N\", Word.InsertLocation.replace);\nawait context.sync();\n\n// A bookmark name must be unique, and the document already carries a\n// \"_GoBack\" bookmark (left over near the end of paragraph 6, marking the\n// previous session's last edit position). Remove it before re-adding\n// \"_GoBack\" at the new edit spot, exactly like Word relocating it itself.\nconst existingGoBack = context.document.body.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nexistingGoBack.load(\"isNullObject\");\nawait context.sync();\nif (!existingGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// Word stamps \"_GoBack\" at the position of the most recent edit, which here\n// is right after \"JUSTIFICACI\u00d3\" (i.e. before the final \"N\" that was typed).\nconst titleParaAfter = context.document.body.paragraphs.items[0];\nconst accentedPart = titleParaAfter.search(\"JUSTIFICACI\u00d3\", { matchCase: true });\naccentedPart.load(\"items\");\nawait context.sync();\n\nconst accentedRange = accentedPart.items[0];\nconst goBackSpot = accentedRange.getRange(Word.RangeLocation.end);\ngoBackSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Change 2: merge the trailing space into the previous sentence --------\n// Paragraph 6 (\"De igual manera ...\") contained that leftover \"_GoBack\"\n// bookmark followed by a run with only a trailing space; now that the\n// bookmark has moved away, fold the trailing space back into the sentence\n// run (single run, same formatting) to match the final text.\nconst paragraphsAfter = context.document.body.paragraphs;\nparagraphsAfter.load(\"text\");\nawait context.sync();\n\nconst lastPara = paragraphsAfter.items[5];\nconst lastParaRange = lastPara.getRange();\nlastParaRange.insertText(\n  \"De igual manera los clientes se ver\u00e1n beneficiados ya que podr\u00e1n consultar los art\u00edculos con los que cuenta la muebler\u00eda desde la comodidad de su hogar y a la hora que lo deseen. \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: title \"JUSTIFICACION\" -> \"JUSTIFICACI\u00d3N\" -------------------\n# Retype the title as \"JUSTIFICACI\u00d3N\" (fixing the missing accent), keep the\n# same run formatting (Arial, bold, size 28).\n$titlePara = $d.Paragraphs(1)\n$titleRange = $titlePara.Range\n$titleRange.Text = \"JUSTIFICACI\u00d3N\"\n\n# Word leaves the \"_GoBack\" bookmark at the position of the last edit, which\n# in this case is right after the inserted \"N\" at the end of the title.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n$titlePara2 = $d.Paragraphs(1)\n$lastEditPos = $titlePara2.Range.End - 1\n$goBackRange = $d.Range($lastEditPos - 1, $lastEditPos - 1)\n$d.Bookmarks.Add(\"_GoBack\", $goBackRange)\n\n# --- Change 2: merge the trailing space into the previous sentence --------\n# Paragraph 6 originally ends in \"...lo deseen.\" + a separate run (after the\n# old \"_GoBack\" bookmark) containing just a trailing space. Merge that\n# trailing space into the sentence run and drop the now-stale bookmark.\n$lastPara = $d.Paragraphs(6)\n$lastParaEnd = $lastPara.Range.End - 1\n$trailingSpaceRange = $d.Range($lastParaEnd - 1, $lastParaEnd)\n$trailingSpaceRange.Text = \"\"\n\n$lastPara2 = $d.Paragraphs(6)\n$insertPoint = $d.Range($lastPara2.Range.End - 1, $lastPara2.Range.End - 1)\n$insertPoint.InsertAfter(\" \")\n"}
